# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the newly scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 4474
$ws1.Range("F13").Value = 686
$ws1.Range("F14").Value = 176
$ws1.Range("F22").Value = 3444
$ws1.Range("F23").Value = 5785
$ws1.Range("F45").Value = 40

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 4474
$ws4.Range("F14").Value = 686
$ws4.Range("F15").Value = 176
$ws4.Range("F23").Value = 3444
$ws4.Range("F24").Value = 5785
$ws4.Range("F42").Value = 895
$ws4.Range("F46").Value = 40
